$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-08 Saturday" "2025-03-09 Sunday"
Replace-Text "35×11=" "13×86="
Replace-Text "94×48=" "93×50="
Replace-Text "47×21=" "38×53="
Replace-Text "58×80=" "16×38="
Replace-Text "29×97=" "36×54="
Replace-Text "92×58=" "69×16="
Replace-Text "75×43=" "42×38="
Replace-Text "50×61=" "50×73="
Replace-Text "80×99=" "91×79="
Replace-Text "51×79=" "69×80="
Replace-Text "96×33=" "48×72="
Replace-Text "64×89=" "92×59="
Replace-Text "57×62=" "75×48="
Replace-Text "84×75=" "41×51="
Replace-Text "85×73=" "45×80="
Replace-Text "54×30=" "50×48="
Replace-Text "96×52=" "43×73="
Replace-Text "41×21=" "27×45="
Replace-Text "53×87=" "42×83="
Replace-Text "36×31=" "83×85="
Replace-Text "39×58=" "24×43="
Replace-Text "88×42=" "39×94="
Replace-Text "91×55=" "27×72="
Replace-Text "63×86=" "74×21="
Replace-Text "26×33=" "81×68="
